$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 63 (shifts existing rows 63-91 down to 64-92),
# inheriting formatting (e.g. the date-style on column D) from the row above.
$ws.Rows.Item(63).Insert()

# Populate the new weekly record in row 63.
$ws.Range("A63").Value = 6
$ws.Range("B63").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C63").Value = "Metropolitana"
$ws.Range("D63").Value = 44879
$ws.Range("E63").Value = 13
$ws.Range("F63").Value = "Fruta"
$ws.Range("G63").Value = 100108
$ws.Range("H63").Value = "Tropicales y subtropicales"
$ws.Range("I63").Value = 100108007
$ws.Range("J63").Value = "Coco"
$ws.Range("K63").Value = "Sin especificar"
$ws.Range("L63").Value = "Primera"
$ws.Range("M63").Value = 150
$ws.Range("N63").Value = 28000
$ws.Range("O63").Value = 30000
$ws.Range("P63").Value = 29000
$ws.Range("Q63").Value = "$/malla 20 unidades"
$ws.Range("R63").Value = "Perú"
$ws.Range("S63").Value = 1450
$ws.Range("T63").Value = 20
